$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.806.98"
$ws.Range("D3").Value = "'2.650.52"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'538.51"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").Value = "'145.75"
$ws.Range("E6").Value = "  +3.46%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.574"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").Value = "'2.665.54"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "'3.127.74"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").Value = "'59.733.58"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "'21.23"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "'2.728.75"
$ws.Range("E17").Value = "  +5.93%  "
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "'345.22"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D21").Value = "'10.39"
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'66.65"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'7.32"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("D29").Value = "'0.0₃0755"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").Value = "'5.85"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'19.04"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").Value = "'150.34"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("E36").Value = "  +1.99%  "
$ws.Range("D37").Value = "'0.846"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "'1.47"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "'293.79"
$ws.Range("E40").Value = "  +5.76%  "
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'0.606"
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("D44").Value = "'19.51"
$ws.Range("E44").Value = "  +5.29%  "
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "'1.983.09"
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").Value = "'18.42"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").Value = "'4.55"
$ws.Range("E51").Value = "  -2.88%  "
